$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values per row (B, C, D columns) to reflect corrected siconfi table filtering

# 2017
$ws.Range("B2").Value = -503.9157522584894
$ws.Range("C2").Value = -184.3024599460749
$ws.Range("D2").Value = -771.4996365605837

# 2018
$ws.Range("B3").Value = -487.7016053425754
$ws.Range("C3").Value = -232.7429042058306
$ws.Range("D3").Value = -514.8661514682541

# 2019
$ws.Range("B4").Value = -516.3149586721969
$ws.Range("C4").Value = -261.4772290195456
$ws.Range("D4").Value = -763.8795634547579

# 2020
$ws.Range("B5").Value = -528.2070474616711
$ws.Range("C5").Value = -316.4916476709561
$ws.Range("D5").Value = -580.4121483922563

# 2021
$ws.Range("B6").Value = -342.8572764143494
$ws.Range("C6").Value = -194.4217271133621

# 2022
$ws.Range("B7").Value = -357.4461201988151
$ws.Range("C7").Value = -215.9024927842444

# 2023
$ws.Range("B8").Value = -370.5427410574931
$ws.Range("C8").Value = -215.6927621177635

# 2024
$ws.Range("B9").Value = -404.1582977226919
$ws.Range("C9").Value = -232.4997573667561
